$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would be re-parsed as a number by COM Value
# auto-typing get NumberFormat forced to Text ("@") first, so they keep
# the sheet's original text/inlineStr representation (e.g. 64.10 stays
# "64.10" instead of being normalized to 64.1).
$ws.Range('D2').Value = '26.146.97'
$ws.Range('E2').Value = '  -4.31%  '
$ws.Range('D3').Value = '1.650.74'
$ws.Range('E3').Value = '  -3.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.23'
$ws.Range('E5').Value = '  -4.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5108'
$ws.Range('E6').Value = '  -3.41%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2580'
$ws.Range('E8').Value = '  -3.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06431'
$ws.Range('E9').Value = '  -3.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.93'
$ws.Range('E10').Value = '  -4.30%  '
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').Value = '1.655.15'
$ws.Range('E12').Value = '  -3.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.278'
$ws.Range('E13').Value = '  -4.95%  '
$ws.Range('D14').Value = '1.879.72'
$ws.Range('E14').Value = '  -3.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5511'
$ws.Range('E15').Value = '  -5.78%  '
$ws.Range('D16').Value = '0.0₅8017'
$ws.Range('E16').Value = '  -2.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.10'
$ws.Range('E17').Value = '  -5.75%  '
$ws.Range('D18').Value = '26.141.90'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '210.54'
$ws.Range('E20').Value = '  -5.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.390'
$ws.Range('E21').Value = '  -5.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.06'
$ws.Range('E22').Value = '  -3.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.050'
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.78'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  +3.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1176'
$ws.Range('E27').Value = '  -2.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.977'
$ws.Range('E28').Value = '  -3.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.78'
$ws.Range('E29').Value = '  -2.75%  '
$ws.Range('E30').Value = '  -3.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.239'
$ws.Range('E31').Value = '  -4.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.342'
$ws.Range('E32').Value = '  -3.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.216'
$ws.Range('E33').Value = '  -6.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.557'
$ws.Range('E34').Value = '  -4.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.737'
$ws.Range('E35').Value = '  -4.64%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.357'
$ws.Range('E36').Value = '  -1.55%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9235'
$ws.Range('E37').Value = '  -3.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5772'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('D39').Value = '1.167.24'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01584'
$ws.Range('E40').Value = '  -3.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.555'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.656'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8265'
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.24'
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('D46').Value = '1.789.67'
$ws.Range('E46').Value = '  -3.47%  '
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4550'
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.46'
$ws.Range('E49').Value = '  -3.82%  '
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.859'
$ws.Range('E51').Value = '  -2.62%  '
